$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.265976
$ws.Range("H2").Value = 36.797928
$ws.Range("I2").Value = 0.004000867643088759
$ws.Range("J2").Value = 0.004000867643088758
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 16.81513620363467
$ws.Range("R2").Value = 151.336225832712
$ws.Range("S2").Value = 0.0000441321030261389
$ws.Range("T2").Value = 0.00004413210302613888
$ws.Range("G3").Value = 12.265976
$ws.Range("H3").Value = 36.797928
$ws.Range("I3").Value = 0.004000867643088759
$ws.Range("J3").Value = 0.004000867643088758
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 1139.677512743283
$ws.Range("R3").Value = 10257.09761468954
$ws.Range("S3").Value = 0.002991136366655686
$ws.Range("T3").Value = 0.002991136366655685
$ws.Range("G4").Value = 12.265976
$ws.Range("H4").Value = 36.797928
$ws.Range("I4").Value = 0.004000867643088759
$ws.Range("J4").Value = 0.004000867643088758
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 364.526211500384
$ws.Range("R4").Value = 3280.735903503456
$ws.Range("S4").Value = 0.0009567159092167034
$ws.Range("T4").Value = 0.0009567159092167029
$ws.Range("G5").Value = 12.265976
$ws.Range("H5").Value = 36.797928
$ws.Range("I5").Value = 0.004000867643088759
$ws.Range("J5").Value = 0.004000867643088758
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 3.384685683416
$ws.Range("R5").Value = 30.462171150744
$ws.Range("S5").Value = 0.00000888326419023145
$ws.Range("T5").Value = 0.000008883264190231447
$ws.Range("I6").Value = 0.9924545876219728
$ws.Range("J6").Value = 0.9924545876219727
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 4171.159997160476
$ws.Range("R6").Value = 37540.43997444429
$ws.Range("S6").Value = 0.01094740241791233
$ws.Range("T6").Value = 0.01094740241791233
$ws.Range("I7").Value = 0.9924545876219728
$ws.Range("J7").Value = 0.9924545876219727
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.7419808086924251
$ws.Range("T7").Value = 0.741980808692425
$ws.Range("I8").Value = 0.9924545876219728
$ws.Range("J8").Value = 0.9924545876219727
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 90424.31372028962
$ws.Range("R8").Value = 813818.8234826066
$ws.Range("S8").Value = 0.2373227954174488
$ws.Range("T8").Value = 0.2373227954174487
$ws.Range("I9").Value = 0.9924545876219728
$ws.Range("J9").Value = 0.9924545876219727
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 839.6045892613646
$ws.Range("R9").Value = 7556.441303352281
$ws.Range("S9").Value = 0.002203581094186575
$ws.Range("T9").Value = 0.002203581094186574
$ws.Range("G10").Value = 8.377189333333332
$ws.Range("H10").Value = 25.131568
$ws.Range("I10").Value = 0.002732438555542716
$ws.Range("J10").Value = 0.002732438555542716
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 11.48409059691911
$ws.Range("R10").Value = 103.356815372272
$ws.Range("S10").Value = 0.00003014052715643161
$ws.Range("T10").Value = 0.00003014052715643161
$ws.Range("G11").Value = 8.377189333333332
$ws.Range("H11").Value = 25.131568
$ws.Range("I11").Value = 0.002732438555542716
$ws.Range("J11").Value = 0.002732438555542716
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 778.3558604054737
$ws.Range("R11").Value = 7005.202743649264
$ws.Range("S11").Value = 0.002042830971240562
$ws.Range("T11").Value = 0.002042830971240562
$ws.Range("G12").Value = 8.377189333333332
$ws.Range("H12").Value = 25.131568
$ws.Range("I12").Value = 0.002732438555542716
$ws.Range("J12").Value = 0.002732438555542716
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 248.9573671676373
$ws.Range("R12").Value = 2240.616304508736
$ws.Range("S12").Value = 0.0006534001297345166
$ws.Range("T12").Value = 0.0006534001297345164
$ws.Range("G13").Value = 8.377189333333332
$ws.Range("H13").Value = 25.131568
$ws.Range("I13").Value = 0.002732438555542716
$ws.Range("J13").Value = 0.002732438555542716
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 2.311610001829333
$ws.Range("R13").Value = 20.804490016464
$ws.Range("S13").Value = 0.000006066927411205504
$ws.Range("T13").Value = 0.000006066927411205504
$ws.Range("G14").Value = 2.489778666666667
$ws.Range("H14").Value = 7.469336
$ws.Range("I14").Value = 0.0008121061793956991
$ws.Range("J14").Value = 0.0008121061793956991
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 3.413178649371555
$ws.Range("R14").Value = 30.718607844344
$ws.Range("S14").Value = 0.000008958045297790903
$ws.Range("T14").Value = 0.000008958045297790902
$ws.Range("G15").Value = 2.489778666666667
$ws.Range("H15").Value = 7.469336
$ws.Range("I15").Value = 0.0008121061793956991
$ws.Range("J15").Value = 0.0008121061793956991
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 231.3346086856809
$ws.Range("R15").Value = 2082.011478171128
$ws.Range("S15").Value = 0.0006071483846691181
$ws.Range("T15").Value = 0.000607148384669118
$ws.Range("G16").Value = 2.489778666666667
$ws.Range("H16").Value = 7.469336
$ws.Range("I16").Value = 0.0008121061793956991
$ws.Range("J16").Value = 0.0008121061793956991
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 73.99244746887466
$ws.Range("R16").Value = 665.932027219872
$ws.Range("S16").Value = 0.0001941966021153434
$ws.Range("T16").Value = 0.0001941966021153433
$ws.Range("G17").Value = 2.489778666666667
$ws.Range("H17").Value = 7.469336
$ws.Range("I17").Value = 0.0008121061793956991
$ws.Range("J17").Value = 0.0008121061793956991
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 0.6870320150586666
$ws.Range("R17").Value = 6.183288135528
$ws.Range("S17").Value = 0.000001803147313446741
$ws.Range("T17").Value = 0.000001803147313446741
